$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values that changed
$ws.Range("F2").Value = "2024-01-30T04:28:37.638579"
$ws.Range("H2").Value = "ACCEPTED_WITH_CONDITIONS"
$ws.Range("J2").Value = "lv.venta.models.users.Student@34c5736d"
$ws.Range("K2").Value = "lv.venta.models.users.AcademicPersonel@4c64a167"
$ws.Range("L2").Value = "[lv.venta.models.users.AcademicPersonel@1001d154]"
$ws.Range("M2").Value = "[lv.venta.models.Comment@1a9accf1, lv.venta.models.Comment@48f3587]"

# Delete row 3 entirely (the second record)
$ws.Rows("3:3").Delete()
